$d = $word.ActiveDocument

function Get-PkgXml($innerParagraphXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
           '<pkg:xmlData>' + `
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + `
           '<w:body>' + $innerParagraphXml + '</w:body></w:document>' + `
           '</pkg:xmlData></pkg:part></pkg:package>'
}

# -----------------------------------------------------------------------
# Change 1: merge the two adjacent runs "dos " + "conforme " (identical
# rPr: sz=28, szCs=28) into a single run "dos conforme " inside the
# paragraph "Sistema verifica validade dos dados conforme DD-visitante."
# (the occurrence that still has "DD-visitante" lower-case / split runs
# — a second, already-merged "DD-Visitante" occurrence also exists and
# must stay untouched, so matching is done case-sensitively).
#
# A plain Find/Replace touching this paragraph coalesces *every*
# identically-formatted adjacent run in the whole paragraph (not only
# the two runs the diff targets), so the paragraph is rebuilt precisely
# via InsertXML instead, preserving every other run exactly as-is.
# -----------------------------------------------------------------------
$found1 = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $t = $para.Range.Text
    if ($t.Contains("verifica validade") -and $t.Contains("DD-visitante")) {
        $rng = $para.Range.Duplicate
        $inner = '<w:p w14:paraId="39D10215" w14:textId="630D15A1" w:rsidR="00C83EB2" w:rsidRDefault="00C83EB2" w:rsidP="00C83EB2">' + `
                 '<w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="18"/></w:numPr><w:ind w:left="360"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr>' + `
                 '<w:r w:rsidRPr="001265C3"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Sistema verifica validade dos da</w:t></w:r>' + `
                 '<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">dos conforme </w:t></w:r>' + `
                 '<w:r w:rsidRPr="001265C3"><w:rPr><w:color w:val="FF0000"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>DD-</w:t></w:r>' + `
                 '<w:r><w:rPr><w:color w:val="FF0000"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>visitante</w:t></w:r>' + `
                 '<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>.</w:t></w:r>' + `
                 '</w:p>'
        $rng.InsertXML((Get-PkgXml $inner))
        $found1 = $true
        break
    }
}
if (-not $found1) {
    throw "Could not locate paragraph for change 1"
}

# -----------------------------------------------------------------------
# Change 2: "Criação da seção Novo Residente" -> split into two runs,
# fixing the spelling/word to "Criação da seção Novo " + "Visitante".
# A plain Find/Replace keeps a single run, so InsertXML is used again to
# get the exact two-run structure shown by the diff.
# -----------------------------------------------------------------------
$found2 = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $t = $para.Range.Text
    if ($t.Contains("seção Novo Residente")) {
        $rng = $para.Range.Duplicate
        $inner = '<w:p w14:paraId="38A31903" w14:textId="19F56780" w:rsidR="006404AA" w:rsidRDefault="006404AA" w:rsidP="004B52A1">' + `
                 '<w:pPr><w:jc w:val="center"/></w:pPr>' + `
                 '<w:r><w:t xml:space="preserve">Criação da seção Novo </w:t></w:r>' + `
                 '<w:r><w:t>Visitante</w:t></w:r>' + `
                 '</w:p>'
        $rng.InsertXML((Get-PkgXml $inner))
        $found2 = $true
        break
    }
}
if (-not $found2) {
    throw "Could not locate paragraph for change 2"
}
